# Updates cryptos list figures (price + 1h volume change) to the latest
# snapshot, and fixes the FirstDigitalUSD / NEARProtocol row ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.474.34"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.191.20"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.22"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.96"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.182.70"
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  +6.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.86"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "3.719.34"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.37"
$ws.Range("E17").Value = "  +4.54%  "
$ws.Range("D18").Value = "3.191.71"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("D19").Value = "64.248.77"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.80"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.742"
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("E24").Value = "  +9.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.22"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +9.51%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.26"
$ws.Range("E32").Value = "  +4.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.119"
$ws.Range("E33").Value = "  +8.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.56"
$ws.Range("E34").Value = "  +7.31%  "
$ws.Range("D35").Value = "0.0₃0867"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.24"
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "468.32"
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("E41").Value = "  +9.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "51.54"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("E43").Value = "  +8.59%  "
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").Value = "2.928.16"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.59"
$ws.Range("E46").Value = "  +8.77%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.14"
$ws.Range("E48").Value = "  +6.75%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  +6.28%  "
$ws.Range("E51").Value = "  +1.70%  "
